# Assignment_Q2_QnA.xlsx edit
# The FAQ sheet gets a new column A: each question is prefixed with
# "FoodFAQ - " and loses its distinct border/fill style (falls back to the
# sheet's default/no style), while column B (answers) keeps its existing
# per-row styling. Column A is also widened slightly and the active
# selection is moved to A8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Question/Answer pairs (row => @(question, answer))
$rows = @(
    @(2,  "FoodFAQ - How can I track my order?", "To track your order, please enter your order number on our tracking page [link]."),
    @(3,  "FoodFAQ - Can I change my order after placing it?", "You can change your order within 24 hours of placing it by contacting our customer support team."),
    @(4,  "FoodFAQ - Can I change my delivery address after placing an order?", "Yes, you can change your delivery address within 24 hours of placing the order by contacting our support team."),
    @(5,  "FoodFAQ - What is the return policy?", "You can return any unopened items within 30 days of purchase. Please visit our returns page for more details."),
    @(6,  "FoodFAQ - Do you offer free shipping?", "Yes, we offer free shipping on orders over `$50."),
    @(7,  "FoodFAQ - Are your products halal-certified?", "Yes, all our products are halal-certified."),
    @(8,  "FoodFAQ - Do you ship internationally?", "Currently, we only ship within Asia."),
    @(9,  "FoodFAQ - How can I apply a discount code?", "You can apply a discount code at checkout. Simply enter the code in the designated field."),
    @(10, "FoodFAQ - What payment methods do you accept?", "We accept credit cards, debit cards, and PayPal."),
    @(11, "FoodFAQ - How do I reset my account password?", "To reset your password, click on 'Forgot Password' on the login page and follow the instructions."),
    @(12, "FoodFAQ - Are there any membership benefits?", "Yes, members receive exclusive discounts, early access to sales, and more. Join our membership program to enjoy these benefits."),
    @(13, "FoodFAQ - Are your products halal-certified?", "Yes, all our products are halal-certified."),
    @(14, "FoodFAQ - Do you offer discounts for bulk purchases?", "Yes, we offer discounts for bulk purchases. Please contact our sales team for more information."),
    @(15, "FoodFAQ - How can I subscribe to your newsletter?", "You can subscribe to our newsletter by entering your email address in the subscription box on our website."),
    @(16, "FoodFAQ - How can I contact customer support?`t", "You can contact our customer support team via email at OnlineSupport@@foodcompany.com or through our online chat.")
)

foreach ($row in $rows) {
    $r = $row[0]
    $question = $row[1]
    $answer = $row[2]

    $ws.Cells.Item($r, 1).Value = $question
    $ws.Cells.Item($r, 1).Style = "Normal"
    $ws.Cells.Item($r, 2).Value = $answer
}

# Column A widened to fit the longer "FoodFAQ - ..." questions.
$ws.Columns("A").ColumnWidth = 53.5

# Move / record the active selection at A8 (matches the saved view state).
$ws.Range("A8").Select()
